$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "HOSP..DE.REHABILITACION.PSICOFISICA"
$ws.Range("B2").Value = "HOSP..DE.REHABILITACION.PSICOFISICA"
$ws.Range("D2").Value = "IREP"
